$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datenerfassung")

# Data for rows 37-47: DateSerial, DurationFraction, LV-Einheit (col C), Bezeichnung (col D)
# NOTE: the rows are listed here in the exact order they were originally authored
# (row 42 was filled in last, after 43-47) so that new shared-string entries get
# appended in the same order as the source workbook.
$rows = @(
    @{ Row=37; DateSerial=44299; Dur=0.10416666666666667; C="Implementierung";                     D="Use Case User Selbstregistrierung" },
    @{ Row=38; DateSerial=44300; Dur=0.125;                C="Implementierung";                     D="Use Case User Selbstregistrierung" },
    @{ Row=39; DateSerial=44300; Dur=0.041666666666666664; C="Koordination und Projektmanagement";  D="Teilnahme an Meeting" },
    @{ Row=40; DateSerial=44301; Dur=0.041666666666666664; C="Implementierung";                     D="Use Case User Selbstregistrierung" },
    @{ Row=41; DateSerial=44305; Dur=0.08333333333333333;  C="Implementierung";                     D="Use Case User Erzeugen durch Admin + Würfeldialog" },
    @{ Row=43; DateSerial=44306; Dur=0.08333333333333333;  C="LV-Einheit";                          D="Teilnahme an Workshop6 Testen" },
    @{ Row=44; DateSerial=44306; Dur=0.0625;                C="Koordination und Projektmanagement"; D="Teilnahme an Meeting" },
    @{ Row=45; DateSerial=44307; Dur=0.08333333333333333;  C="Implementierung";                     D="User-related Tests" },
    @{ Row=46; DateSerial=44308; Dur=0.10416666666666667;  C="Implementierung";                     D="User-related Tests" },
    @{ Row=47; DateSerial=44309; Dur=0.125;                C="Implementierung";                     D="User-related Tests" },
    @{ Row=42; DateSerial=44306; Dur=0.08333333333333333;  C="Implementierung";                     D="Use Case User Editieren durch Admin" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.DateSerial
    $ws.Cells.Item($r.Row, 2).Value = $r.Dur
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $cellD = $ws.Cells.Item($r.Row, 4)
    $cellD.Value = $r.D
    $cellD.NumberFormat = "@"
    $cellD.Font.Italic = $true
    $cellD.HorizontalAlignment = 1
}

# Update view state: scroll the frozen (bottom-right) pane down and set the
# active selection to match the edited area.
$win = $excel.ActiveWindow
$win.ScrollRow = 20
$win.ScrollColumn = 3
$ws.Range("D49").Select()
